$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Step A: Drop the "Meta description: ..." paragraph. Its trailing
# "Gameplay Mechanics" run re-joins the Heading1 title paragraph, i.e.
# the paragraph break + "Meta description: ...now." text disappear and
# "...Review" is immediately followed by "Gameplay Mechanics" under the
# same Heading1-styled paragraph.
# -------------------------------------------------------------------

$titleRange = $d.Content
$titleRange.Find.Execute("Play Ghostbusters Plus Free Slot Game - Review", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$titleEnd = $titleRange.End

$gmRange = $d.Content
$gmRange.Find.Execute("Gameplay Mechanics", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$gmStart = $gmRange.Start

# Remove the "Meta description: ...now." text (between the end of the
# title text and the start of "Gameplay Mechanics"), keeping the
# paragraph mark itself for now.
$d.Range($titleEnd + 1, $gmStart).Delete()

# Collapse the now-orphaned paragraph mark so the title paragraph and
# the (now empty) former meta-description paragraph become one.
$markRange = $d.Range($titleEnd, $titleEnd)
$markRange.Collapse(1)
$markRange.Delete(1, 1)

# A paragraph-mark merge adopts the *following* paragraph's properties,
# so re-apply Heading1 explicitly to the surviving (first) paragraph.
$d.Paragraphs.Item(1).Style = "Heading1"

# -------------------------------------------------------------------
# Step B: Insert a new bold paragraph - "Play Ghostbusters Plus Free
# Slot Game - Review" - immediately before the closing ("Prompt: ...")
# paragraph.
# -------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$d.Paragraphs.Item($lastIndex - 1).Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($lastIndex)
$newPara.Style = "Normal"

$newParaStart = $newPara.Range.Start
$insertPoint = $d.Range($newParaStart, $newParaStart)
$insertPoint.InsertAfter("Play Ghostbusters Plus Free Slot Game - Review")
$insertPoint.Font.Bold = $true

# -------------------------------------------------------------------
# Step C: Swap the old "Prompt: ..." image-prompt text (now the very
# last paragraph) for the meta-description copy, keeping the run's
# italic formatting intact.
# -------------------------------------------------------------------

$oldPrompt = "Prompt: Please create a feature image in cartoon style for the game " + [char]34 + "Ghostbuster Plus" + [char]34 + ". The image should include a happy Maya warrior with glasses. The Maya warrior should be holding a ghost-catching tool and standing in front of the Ghostbusters headquarters. The background of the image should feature a cityscape with ghosts flying around. The colors should be bright and playful to match the tone of the game. The size of the image should be 1080x1080 pixels."
$newMeta = "Read our review of Ghostbusters Plus slot game, a cinema-themed online slot that offers bonus features. Play the free demo now."

$promptRange = $d.Content
$promptRange.Find.ClearFormatting()
$promptRange.Find.Execute($oldPrompt, $false, $false, $false, $false, $false, $true, 1, $false, $newMeta, 2) | Out-Null
